$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1983
$ws.Range("I18").Value = 1983
$ws.Range("K18").Value = 1983
$ws.Range("M18").Value = -1699

$ws.Range("H19").Value = 1380.762
$ws.Range("I19").Value = 1238.4615
$ws.Range("K19").Value = 1238.4615
$ws.Range("M19").Value = -1063.4615

$ws.Range("H40").Value = 4151.609
$ws.Range("J40").Value = 4540.6665
$ws.Range("L40").Value = 4540.6665
$ws.Range("N40").Value = -4890.6665

$ws.Range("H52").Value = 1000
$ws.Range("I52").Value = 1000
$ws.Range("K52").Value = 3000
$ws.Range("M52").Value = -2840

$ws.Range("H58").Value = 325.8
$ws.Range("I58").Value = 195.44444
$ws.Range("K58").Value = 586.33332
$ws.Range("M58").Value = -436.33332

$ws.Range("H98").Value = 1207.1945
$ws.Range("I98").Value = 955.97144
$ws.Range("K98").Value = 955.97144
$ws.Range("M98").Value = 542.02856

$ws.Range("H107").Value = 690.75
$ws.Range("I107").Value = 756.2222
$ws.Range("J107").Value = 494.33334
$ws.Range("K107").Value = 756.2222
$ws.Range("L107").Value = 494.33334
$ws.Range("M107").Value = 1163.7778
$ws.Range("N107").Value = -4334.33334

$ws.Range("H122").Value = 1207.1945
$ws.Range("I122").Value = 955.97144
$ws.Range("K122").Value = 2867.91432
$ws.Range("M122").Value = -417.9143199999999

$ws.Range("H123").Value = 300000.5
$ws.Range("J123").Value = 300000.5
$ws.Range("L123").Value = 300000.5
$ws.Range("N123").Value = -309800.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 104
$ws.Range("I19").Value = 104
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 104
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 125
$ws.Range("N19").ClearContents()

$ws.Range("H61").Value = 50002084
$ws.Range("I61").Value = 52633640
$ws.Range("K61").Value = 52633640
$ws.Range("M61").Value = -52633428

$ws.Range("H74").Value = 50003548
$ws.Range("I74").Value = 50003548
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 50003548
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -50002674
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 50003548
$ws.Range("I77").Value = 50003548
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 250017740
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -250013372
$ws.Range("N77").ClearContents()

$ws.Range("H102").Value = 11112854
$ws.Range("I102").Value = 14287399
$ws.Range("K102").Value = 14287399
$ws.Range("M102").Value = -14285777

$ws.Range("H122").Value = 4668.303
$ws.Range("I122").Value = 3440.7693
$ws.Range("K122").Value = 10322.3079
$ws.Range("M122").Value = -7872.3079

$ws.Range("H132").Value = 3228604.5
$ws.Range("I132").Value = 3336041.5
$ws.Range("K132").Value = 10008124.5
$ws.Range("M132").Value = -10005594.5

$ws.Range("H136").Value = 50002084
$ws.Range("I136").Value = 52633640
$ws.Range("K136").Value = 157900920
$ws.Range("M136").Value = -157898370

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1008.7778
$ws.Range("I20").Value = 1219.3334
$ws.Range("K20").Value = 1219.3334
$ws.Range("M20").Value = -972.3334

$ws.Range("H134").Value = 15629579
$ws.Range("I134").Value = 16133605
$ws.Range("K134").Value = 48400815
$ws.Range("M134").Value = -48398280

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 22599
$ws.Range("I22").Value = 50499
$ws.Range("J22").Value = 3999
$ws.Range("K22").Value = 50499
$ws.Range("L22").Value = 3999
$ws.Range("M22").Value = -50149
$ws.Range("N22").Value = -4699

$ws.Range("H31").Value = 7451
$ws.Range("I31").Value = 4938.625
$ws.Range("K31").Value = 4938.625
$ws.Range("M31").Value = -4643.625

$ws.Range("H34").Value = 7451
$ws.Range("I34").Value = 4938.625
$ws.Range("K34").Value = 4938.625
$ws.Range("M34").Value = -4736.625

$ws.Range("H39").Value = 24816.666
$ws.Range("I39").Value = 4451
$ws.Range("K39").Value = 4451
$ws.Range("M39").Value = -4060

$ws.Range("H49").Value = 24816.666
$ws.Range("I49").Value = 4451
$ws.Range("K49").Value = 4451
$ws.Range("M49").Value = -4269

$ws.Range("H54").Value = 33199.6
$ws.Range("J54").Value = 33199.6
$ws.Range("L54").Value = 33199.6
$ws.Range("N54").Value = -34515.6

$ws.Range("H132").Value = 200002900
$ws.Range("I132").Value = 333335500
$ws.Range("K132").Value = 1000006500
$ws.Range("M132").Value = -1000003970

$ws.Range("H134").Value = 19233912
$ws.Range("I134").Value = 20836414
$ws.Range("K134").Value = 62509242
$ws.Range("M134").Value = -62506707

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 362.33334
$ws.Range("I98").Value = 362.33334
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1087.00002
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 410.9999800000001
$ws.Range("N98").ClearContents()

$ws.Range("H114").Value = 126236.625
$ws.Range("I114").Value = 166982.33
$ws.Range("J114").Value = 3999.5
$ws.Range("K114").Value = 500946.99
$ws.Range("L114").Value = 11998.5
$ws.Range("M114").Value = -497692.99
$ws.Range("N114").Value = -18506.5

$ws.Range("H115").Value = 5749.5
$ws.Range("I115").Value = 4999
$ws.Range("K115").Value = 14997
$ws.Range("M115").Value = -13822

$ws.Range("H129").Value = 4987
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 4987
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 14961
$ws.Range("N129").Value = -24961
$ws.Range("M129").ClearContents()

$ws.Range("H131").Value = 1669.75
$ws.Range("J131").Value = 2198.6
$ws.Range("L131").Value = 6595.799999999999
$ws.Range("N131").Value = -16675.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3593.6875
$ws.Range("I102").Value = 3593.6875
$ws.Range("K102").Value = 3593.6875
$ws.Range("M102").Value = -1971.6875

$ws.Range("H122").Value = 3376.6
$ws.Range("I122").Value = 1642.1818
$ws.Range("J122").Value = 8146.25
$ws.Range("K122").Value = 4926.5454
$ws.Range("L122").Value = 24438.75
$ws.Range("M122").Value = -2476.5454
$ws.Range("N122").Value = -29338.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2442.0454
$ws.Range("I22").Value = 1982.4615
$ws.Range("J22").Value = 3105.889
$ws.Range("K22").Value = 1982.4615
$ws.Range("L22").Value = 3105.889
$ws.Range("M22").Value = -1687.4615
$ws.Range("N22").Value = -3695.889

$ws.Range("H27").Value = 2442.0454
$ws.Range("I27").Value = 1982.4615
$ws.Range("J27").Value = 3105.889
$ws.Range("K27").Value = 1982.4615
$ws.Range("L27").Value = 3105.889
$ws.Range("M27").Value = -1875.4615
$ws.Range("N27").Value = -3319.889

$ws.Range("H40").Value = 4588.4443
$ws.Range("I40").Value = 4588.4443
$ws.Range("K40").Value = 4588.4443
$ws.Range("M40").Value = -4452.4443

$ws.Range("H61").Value = 1361.5714
$ws.Range("I61").Value = 1256.2
$ws.Range("J61").Value = 1625
$ws.Range("K61").Value = 1256.2
$ws.Range("L61").Value = 1625
$ws.Range("M61").Value = -1054.2
$ws.Range("N61").Value = -2029

$ws.Range("H113").Value = 1361.5714
$ws.Range("I113").Value = 1256.2
$ws.Range("J113").Value = 1625
$ws.Range("K113").Value = 1256.2
$ws.Range("L113").Value = 1625
$ws.Range("M113").Value = 913.8
$ws.Range("N113").Value = -5965

$ws.Range("H136").Value = 1279.1177
$ws.Range("I136").Value = 1135.1482
$ws.Range("K136").Value = 3405.4446
$ws.Range("M136").Value = -855.4446000000003

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 13159759
$ws.Range("I136").Value = 13159759
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 39479277
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -39476727
$ws.Range("N136").ClearContents()
